$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 66, shifting rows 66:277 down to 67:278
$ws.Rows("66:66").Insert()

# Fill the new row 66 with the new data values (matching the standard pattern
# used throughout this sheet), with the new Fecha (D) and Volumen (J) values.
$ws.Range("A66").Value = 3
$ws.Range("B66").Value = "Femacal de La Calera"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44623
$ws.Range("E66").Value = 5
$ws.Range("F66").Value = 100112039
$ws.Range("G66").Value = "Ciboulette"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 160
$ws.Range("K66").Value = 1500
$ws.Range("L66").Value = 1500
$ws.Range("M66").Value = 1500
$ws.Range("N66").Value = "$/docena de atados"
$ws.Range("O66").Value = "Provincia de Quillota"
$ws.Range("P66").Value = 500
$ws.Range("Q66").Value = 3
$ws.Range("R66").Value = "Hortaliza"
